$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C (_circuit) and I (patchpanelport) hold numeric-looking text
# ("1", "2", "4") in the source data, so force text storage before writing.
$ws.Range("C2:C4").NumberFormat = "@"
$ws.Range("I2:I4").NumberFormat = "@"

# --- Update existing row 2 (site renamed iad55/iad -> gru1/gru, provider akamai -> google) ---
$ws.Range("A2").Value = "gru1"
$ws.Range("B2").Value = "gru"
$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "google"
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = "gru1-br-cor-r3"
$ws.Range("G2").Value = "et-0/0/11"
$ws.Range("H2").Value = "pp1"
$ws.Range("I2").Value = "1"

# --- New row 3 (circuit 2, microsoft, gru1-br-tra-r3) ---
$ws.Range("A3").Value = "gru1"
$ws.Range("B3").Value = "gru"
$ws.Range("C3").Value = "2"
$ws.Range("D3").Value = "microsoft"
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = "gru1-br-tra-r3"
$ws.Range("G3").Value = "et-0/0/11"
$ws.Range("H3").Value = "pp1"
$ws.Range("I3").Value = "2"

# --- New row 4 (circuit 4, globenet, gru1-br-cor-r4) ---
$ws.Range("A4").Value = "gru1"
$ws.Range("B4").Value = "gru"
$ws.Range("C4").Value = "4"
$ws.Range("D4").Value = "globenet"
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = "gru1-br-cor-r4"
$ws.Range("G4").Value = "et-0/0/37"
$ws.Range("H4").Value = "pp1"
$ws.Range("I4").Value = "4"
